# Performance Food order sheet — add two new line items (rows 50 and 51)
# to the bottom of the existing order table, matching the existing
# "everything stored as text" convention used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: Mushroom Portobello - 4" Cap
$ws.Range("A50").Value = "'20478"
$ws.Range("B50").Value = "Mushroom Portobello - 4`" Cap"
$ws.Range("C50").Value = "'6"
$ws.Range("D50").Value = "'21.80"
$ws.Range("E50").Value = "'130.80"

# Row 51: Strawberry - Fresh
$ws.Range("A51").Value = "'39142"
$ws.Range("B51").Value = "Strawberry - Fresh"
$ws.Range("C51").Value = "'6"
$ws.Range("D51").Value = "'30.80"
$ws.Range("E51").Value = "'184.80"
